$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data row 2 with new generic test data
$ws.Range("A2").Value = "First"
$ws.Range("B2").Value = "Member"
$ws.Range("C2").Value = "first_member@mail.com"

# Column A first (Second/Third/Fourth), then B, then column C for rows 3-5
$ws.Range("A3").Value = "Second"
$ws.Range("A4").Value = "Third"
$ws.Range("A5").Value = "Fourth"

$ws.Range("B3").Value = "Member"
$ws.Range("B4").Value = "Member"
$ws.Range("B5").Value = "Member"

$ws.Range("C3").Value = "second_member@mail.com"
$ws.Range("C4").Value = "third_member@mail.com"
$ws.Range("C5").Value = "fourth_member@mail.com"

$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:third_member@mail.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:fourth_member@mail.com")

$ws.Range("C4").Style = "Hyperlink"
$ws.Range("C5").Style = "Hyperlink"

$wb.Save()
